# Insert a new weekly price record at row 107 for "Poroto granado" (Femacal de
# La Calera), pushing the existing rows 107-205 down to 108-206.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 107, shifting rows 107:205 -> 108:206.
$ws.Rows.Item(107).Insert()

# Populate the new row 107 with the new record's data.
$ws.Range("A107").Value = 3
$ws.Range("B107").Value = "Femacal de La Calera"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = 44904
$ws.Range("E107").Value = 5
$ws.Range("F107").Value = 100112030
$ws.Range("G107").Value = "Poroto granado"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 38
$ws.Range("K107").Value = 41000
$ws.Range("L107").Value = 41000
$ws.Range("M107").Value = 41000
$ws.Range("N107").Value = "`$/saco 25 kilos"
$ws.Range("O107").Value = "Provincia de Limarí"
$ws.Range("P107").Value = 1640
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"
